$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect against Excel auto-converting numeric-looking strings to numbers/dates
# by forcing text format before assignment, then reset cell style to Normal
# afterwards so no visual/style diff is introduced.
$targetCells = @('D2', 'E2', 'D3', 'E3', 'E4', 'D5', 'E5', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E10', 'D11', 'E11', 'E12', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'D18', 'E18', 'D19', 'E19', 'D20', 'E21', 'D22', 'E22', 'D23', 'E23', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'D28', 'E28', 'D29', 'E29', 'E30', 'D31', 'E31', 'D32', 'E32', 'D33', 'E33', 'D34', 'E34', 'D35', 'E35', 'B36', 'C36', 'D36', 'E36', 'B37', 'C37', 'D37', 'E37', 'B38', 'C38', 'D38', 'E38', 'D39', 'E39', 'B40', 'C40', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45', 'E46', 'D47', 'E47', 'D48', 'E48', 'D49', 'E49', 'D50', 'E50', 'D51', 'E51')
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '27.975.27'
$ws.Range('E2').Value = '  +1.52%  '
$ws.Range('D3').Value = '1.753.46'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').Value = '336.01'
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('E6').Value = '  -0.43%  '
$ws.Range('D7').Value = '0.3832'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '0.3428'
$ws.Range('E8').Value = '  +0.70%  '
$ws.Range('D9').Value = '46.04'
$ws.Range('E9').Value = '  -1.64%  '
$ws.Range('D10').Value = '1.122'
$ws.Range('E10').Value = '  -1.27%  '
$ws.Range('D11').Value = '0.07230'
$ws.Range('E11').Value = '  -1.99%  '
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('D14').Value = '6.169'
$ws.Range('E14').Value = '  -2.60%  '
$ws.Range('D15').Value = '7.147'
$ws.Range('E15').Value = '  +1.81%  '
$ws.Range('D16').Value = '1.750.38'
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('D17').Value = '0.00001061'
$ws.Range('E17').Value = '  -1.15%  '
$ws.Range('D18').Value = '0.06599'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('D19').Value = '79.31'
$ws.Range('E19').Value = '  -3.30%  '
$ws.Range('D20').Value = '0.9991'
$ws.Range('E21').Value = '  -3.45%  '
$ws.Range('D22').Value = '6.199'
$ws.Range('E22').Value = '  -3.03%  '
$ws.Range('D23').Value = '27.986.45'
$ws.Range('E23').Value = '  +1.57%  '
$ws.Range('E24').Value = '  -3.15%  '
$ws.Range('D25').Value = '2.381'
$ws.Range('E25').Value = '  +0.61%  '
$ws.Range('D26').Value = '154.04'
$ws.Range('E26').Value = '  +0.64%  '
$ws.Range('D27').Value = '19.84'
$ws.Range('E27').Value = '  -3.86%  '
$ws.Range('D28').Value = '2.300'
$ws.Range('E28').Value = '  -5.46%  '
$ws.Range('D29').Value = '1.951.62'
$ws.Range('E29').Value = '  -0.56%  '
$ws.Range('E30').Value = '  -10.99%  '
$ws.Range('D31').Value = '131.36'
$ws.Range('E31').Value = '  -2.40%  '
$ws.Range('D32').Value = '4.024'
$ws.Range('E32').Value = '  +1.54%  '
$ws.Range('D33').Value = '5.842'
$ws.Range('E33').Value = '  -4.11%  '
$ws.Range('D34').Value = '0.08812'
$ws.Range('E34').Value = '  +0.36%  '
$ws.Range('D35').Value = '12.22'
$ws.Range('E35').Value = '  -3.77%  '
$ws.Range('B36').Value = 'WEMIXTOKEN'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '1.541'
$ws.Range('E36').Value = '  +2.85%  '
$ws.Range('B37').Value = 'TheSandbox'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D37').Value = '0.6571'
$ws.Range('E37').Value = '  -3.00%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.02288'
$ws.Range('E38').Value = '  -4.86%  '
$ws.Range('D39').Value = '5.152'
$ws.Range('E39').Value = '  -3.47%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '0.06163'
$ws.Range('E40').Value = '  -1.64%  '
$ws.Range('D41').Value = '0.2104'
$ws.Range('E41').Value = '  -3.45%  '
$ws.Range('D42').Value = '1.214'
$ws.Range('E42').Value = '  -2.66%  '
$ws.Range('D43').Value = '7.953'
$ws.Range('E43').Value = '  -3.84%  '
$ws.Range('D44').Value = '0.9982'
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('D45').Value = '13.73'
$ws.Range('E45').Value = '  -3.08%  '
$ws.Range('E46').Value = '  +0.30%  '
$ws.Range('D47').Value = '0.6061'
$ws.Range('E47').Value = '  -2.76%  '
$ws.Range('D48').Value = '127.10'
$ws.Range('E48').Value = '  -3.52%  '
$ws.Range('D49').Value = '2.009'
$ws.Range('E49').Value = '  -3.01%  '
$ws.Range('D50').Value = '1.168'
$ws.Range('E50').Value = '  +2.11%  '
$ws.Range('D51').Value = '1.111'
$ws.Range('E51').Value = '  +4.64%  '

foreach ($addr in $targetCells) {
    $ws.Range($addr).Style = "Normal"
}
